$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: extend the description text for the button-head-screws note ---
$ws.Range("D2").Value = "Button head screws won't work due to its head size, Hex is preferred for its resilience against stripping"

# --- Row 3: remove the now-merged "Hex is preferred..." note, right-align the qty ---
$ws.Range("D3").ClearContents()
$ws.Range("C3").HorizontalAlignment = -4152  # xlRight -> creates style index 3

# --- Row 5: right-align qty, update the insert description ---
$ws.Range("C5").HorizontalAlignment = -4152  # reuse style index 3
$ws.Range("D5").Value = "18 is needed if extension is used, M3 insert, 4mm OD, 5mm Length (shorter length could still work)"

# --- Row 7: repurposed as the new "Extensions extra optional" note ---
$ws.Range("A7").Value = "Extensions extra optional"
$ws.Range("D7").Value = "Add these numbers up for each extension that you are using"

# --- Row 8: repurposed as a new BOM line (M3x6 BHCS for base extension) ---
$ws.Range("A8").Value = "M3x6 BHCS"
$ws.Range("B2").Copy()
$ws.Range("B8").PasteSpecial(-4122)  # xlPasteFormats, reuse style index 1
$ws.Range("B8").Value = "ISO 4762 / DIN 912"
$ws.Range("C8").Value = 4
$ws.Range("D8").Value = "Base extension screws"

# --- Row 9: repurposed as a new BOM line (M3x8 BHCS for ramp extension) ---
$ws.Range("A9").Value = "M3x8 BHCS"
$ws.Range("B2").Copy()
$ws.Range("B9").PasteSpecial(-4122)  # xlPasteFormats, reuse style index 1
$ws.Range("B9").Value = "ISO 4762 / DIN 912"
$ws.Range("C9").Value = 4
$ws.Range("D9").Value = "Ramp extension screws"

# --- Row 10 (new): base-extension heatset-insert quantity ---
$ws.Range("A10").Value = "M3x4x5 Brass Heatset Insert"
$ws.Range("B10").WrapText = $true            # creates style index 4
$ws.Range("C10").Value = 4
$ws.Range("C10").HorizontalAlignment = -4152 # reuse style index 3
$ws.Range("A10").HorizontalAlignment = -4131 # xlLeft -> creates style index 5
$ws.Range("D10").Value = "4 is needed for a base extension"

# --- Row 11 (new): ramp-extension heatset-insert quantity, merged with A10 ---
$ws.Range("B11").WrapText = $true            # reuse style index 4
$ws.Range("C11").Value = "4 or 6"
$ws.Range("C11").HorizontalAlignment = -4152 # reuse style index 3
$ws.Range("A11").HorizontalAlignment = -4131 # reuse style index 5
$ws.Range("D11").Value = "4 is needed for a ramp extension if you are using a base extension, if you aren't then you need 6"
$ws.Range("D11").Font.Color = 0              # creates style index 6 (black Aptos Narrow)

$ws.Range("A10:A11").Merge()

# --- Rows 13-15: the "Other tools needed" block, moved down from rows 7-9 ---
$ws.Range("A13").Value = "Other tools needed"
$ws.Range("A14").Value = "Soldering iron"
$ws.Range("A15").Value = "Screwdriver with a hex head/hex key"

# --- Column D is now much wider to fit the longer descriptions ---
$ws.Columns.Item(4).ColumnWidth = 88.6666666666667

# --- Selection moves to D13 ---
$ws.Range("D13").Select()

# --- Header/footer classification banner ---
$ws.PageSetup.CenterHeader = "&""Calibri""&10&KFFFF00 RMIT Classification: Trusted&1#`r"
